# Insert a new data row at row 131, shifting existing rows 131:194 down to 132:195.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("131").Insert()

# Populate the newly inserted row 131 with the new record's data (same
# constant columns as the surrounding rows, plus the new observation values).
$ws.Range("A131").Value = 4
$ws.Range("B131").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C131").Value = "Los Lagos"
$ws.Range("D131").Value = 44553
$ws.Range("E131").Value = 10
$ws.Range("F131").Value = "Fruta"
$ws.Range("G131").Value = 100102
$ws.Range("H131").Value = "Cítricos"
$ws.Range("I131").Value = 100102006
$ws.Range("J131").Value = "Pomelo"
$ws.Range("K131").Value = "Start Ruby"
$ws.Range("L131").Value = "Primera"
$ws.Range("M131").Value = 200
$ws.Range("N131").Value = 10000
$ws.Range("O131").Value = 11000
$ws.Range("P131").Value = 10500
$ws.Range("Q131").Value = "$/caja 14 kilos empedrada"
$ws.Range("R131").Value = "Región de O'Higgins"
$ws.Range("S131").Value = 750
$ws.Range("T131").Value = 14
